$wb = $excel.ActiveWorkbook

# Sheet "展览": F4 1227 -> 1231, F5 613 -> 616
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1231
$wsExhibit.Range("F5").Value = 616

# Sheet "全部类型": F4 1227 -> 1231, F6 613 -> 616
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1231
$wsAll.Range("F6").Value = 616
